$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Get-FindRange($searchText) {
    $rng = $d.Content.Duplicate
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng
}

# Adds "sudo " (Times font) before $afterText, removing a leading "# " if present right before it.
function Add-SudoPrefix($afterText) {
    $full = "# " + $afterText
    $rng = Get-FindRange($full)
    $hashRange = $d.Range($rng.Start, $rng.Start + 2)
    $hashRange.Text = ""
    $insertPoint = $d.Range($rng.Start, $rng.Start)
    $insertPoint.InsertBefore("sudo ")
    $sudoRange = $d.Range($rng.Start, $rng.Start + 5)
    $sudoRange.Font.Name = "Times"
}

# ---------------------------------------------------------------------------
# 1. Fix hyphenation artifact "de- signs" -> "designs"
Replace-Text "protocol de- signs and implementations" "protocol designs and implementations"

# ---------------------------------------------------------------------------
# 2. "# sysctl -q net.ipv4.tcp_max_syn_backlog" -> "sudo " + bookmark + "sysctl -q ..."
Add-SudoPrefix "sysctl -q net.ipv4.tcp_max_syn_backlog"
$rng2 = Get-FindRange("sysctl -q net.ipv4.tcp_max_syn_backlog")
$bmRange = $d.Range($rng2.Start, $rng2.Start)
$d.Bookmarks.Add("__DdeLink__141_514673332", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. Fix "packaet" typo + insert "to that IP address"
Replace-Text "to send one packaet at a time, and" "to send one packet at a time to that IP address, and"

# ---------------------------------------------------------------------------
# 4. "...attack is successful or not.  " -> "...attack has potential to succeed."
# (done further below, AFTER inserting the three new paragraphs, so the
#  search text used for positioning the new paragraphs stays unambiguous)

# ---------------------------------------------------------------------------
# 5. Insert three new paragraphs about shrinking the backlog queue
$rng = Get-FindRange("attack is successful or not.  ")
$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertParagraphAfter()
$p1 = $d.Range($rng.End + 1, $rng.End + 1)
$p1.InsertAfter("To make your attack easier to succeed, we will shrink the size of the backlog queue to 5:")

$rngP1 = Get-FindRange("backlog queue to 5:")
$insertPoint2 = $d.Range($rngP1.End, $rngP1.End)
$insertPoint2.InsertParagraphAfter()
$p2 = $d.Range($rngP1.End + 1, $rngP1.End + 1)
$p2.InsertAfter("sudo sysctl -w net.ipv4.tcp_max_syn_backlog=5")
$p2para = $p2.Paragraphs.First
$p2para.Format.FirstLineIndent = 17.05

$rngP2 = Get-FindRange("tcp_max_syn_backlog=5")
$insertPoint3 = $d.Range($rngP2.End, $rngP2.End)
$insertPoint3.InsertParagraphAfter()
$p3 = $d.Range($rngP2.End + 1, $rngP2.End + 1)
$p3.InsertAfter("Send five packets via nping and then try to telnet to the server via the user component.  Report on your success.")
$p3para = $p3.Paragraphs.First
$p3para.Format.FirstLineIndent = -0.01

# Now perform replacement #4
Replace-Text "describe how you know whether the attack is successful or not.  " "describe how you know whether the attack has potential to succeed."

# ---------------------------------------------------------------------------
# 9. "Relative Sequence Number and Window Scaling" -> "Relative Sequence Number"
Replace-Text "Relative Sequence Number and Window Scaling" "Relative Sequence Number"

# ---------------------------------------------------------------------------
# 6,7,8. sudo prefixes on sysctl commands
Add-SudoPrefix "sysctl -a | grep cookie (Display the SYN cookie flag) "
Add-SudoPrefix "sysctl -w net.ipv4.tcp_syncookies=0 (turn off SYN cookie) "
Add-SudoPrefix "sysctl -w net.ipv4.tcp_syncookies=1 (turn on SYN cookie)"

# ---------------------------------------------------------------------------
# 10. "3.4 Task 4 : TCP Session Hijacking" -> "3.3 Task 3 : TCP Session Hijacking"
Replace-Text "3.4 Task 4 : TCP Session Hijacking" "3.3 Task 3 : TCP Session Hijacking"

# ---------------------------------------------------------------------------
# 11. typo "npig" -> "nping"
Replace-Text "packet spoofing (npig) to perform" "packet spoofing (nping) to perform"

# ---------------------------------------------------------------------------
# 12. "3.5 Task 5 : Creating Reverse Shell using TCP Session Hijacking" -> "3.4 Task 4 : ..."
Replace-Text "3.5 Task 5 : Creating Reverse Shell using TCP Session Hijacking" "3.4 Task 4 : Creating Reverse Shell using TCP Session Hijacking"

# ---------------------------------------------------------------------------
# 14. The empty paragraph right after "moreterm.py tcpip attacker" gains
#     explicit run formatting (Times/11pt/black) matching its own pPr/rPr.
$rngMore = Get-FindRange("moreterm.py tcpip attacker")
$moreParaStart = $rngMore.Paragraphs.First.Range.Start
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $moreParaStart) {
        $nextP = $d.Paragraphs.Item($i + 1)
        $nextP.Range.Font.Name = "Times"
        $nextP.Range.Font.Size = 11
        $nextP.Range.Font.SizeBi = 11
        $nextP.Range.Font.Color = 0
        break
    }
}
